# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 07:35"

# --- Row 70: Hungria - numeric refresh only ---
$ws.Range("B70").Value = 3598
$ws.Range("C70").Value = 42
$ws.Range("D70").Value = 1454
$ws.Range("E70").Value = 1674
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 470

# --- Row 73: Tailandia - numeric refresh only ---
$ws.Range("B73").Value = 3034
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = 2888
$ws.Range("E73").Value = 90
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 56

# --- Rows 75-76: Guinea / Uzbekistan swap ranking, Uzbekistan data refreshed ---
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 2880
$ws.Range("C75").Value = 25
$ws.Range("D75").Value = 2338
$ws.Range("E75").Value = 529
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 13

$ws.Range("A76").Value = "Guinea"
$ws.Range("B76").Value = 2863
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 1525
$ws.Range("E76").Value = 1320
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 18

# --- Row 81: Bulgaria - numeric refresh only ---
$ws.Range("B81").Value = 2292
$ws.Range("C81").Value = 33
$ws.Range("D81").Value = 684
$ws.Range("E81").Value = 1492
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 116

# --- Row 96: El Salvador - numeric refresh only ---
$ws.Range("E96").Value = 965
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 31

# --- Rows 195-196: Santa Lucia / Nueva Caledonia swap ranking (data identical) ---
$ws.Range("A195").Value = "Nueva Caledonia"
$ws.Range("A196").Value = "Santa Lucia"

# --- Rows 209, 211: Seychelles / Montserrat swap ranking, data refreshed ---
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# --- Rows 214-215: Sahara Occidental / San Bartolome swap ranking (data identical) ---
$ws.Range("A214").Value = "San Bartolome"
$ws.Range("A215").Value = "Sahara Occidental"
